$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (Coinranking snapshot) on Fri Dec 22 04:28:00 UTC 2023 with GitHub Actions.
# Refreshes Price (col D) and Volume(1h) (col E) for each coin row; rows 41/42 additionally
# swap their Coin/Link/Price/Volume content (Celestia <-> LidoDAOToken reordering).
# Numeric-looking price strings are written with a leading apostrophe so Excel keeps them
# as text (matching the source data, which stores prices as text, not numbers).

$ws.Range("D2").Value = "44.091.45"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "2.251.94"
$ws.Range("E3").Value = "  +2.28%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'98.69"
$ws.Range("E5").Value = "  +17.05%  "
$ws.Range("D6").Value = "'272.31"
$ws.Range("E6").Value = "  +5.53%  "
$ws.Range("E7").Value = "  +1.76%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +4.93%  "
$ws.Range("D10").Value = "'48.14"
$ws.Range("E10").Value = "  +7.74%  "
$ws.Range("D11").Value = "'0.0942"
$ws.Range("E11").Value = "  +2.18%  "
$ws.Range("E12").Value = "  +14.18%  "
$ws.Range("E13").Value = "  -0.07%  "
$ws.Range("D14").Value = "'15.41"
$ws.Range("E14").Value = "  +7.68%  "
$ws.Range("D15").Value = "2.582.39"
$ws.Range("E15").Value = "  +2.07%  "
$ws.Range("D16").Value = "'0.826"
$ws.Range("E16").Value = "  +5.86%  "
$ws.Range("D17").Value = "2.256.81"
$ws.Range("E17").Value = "  +2.56%  "
$ws.Range("D18").Value = "44.089.42"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("E19").Value = "  +2.75%  "
$ws.Range("D20").Value = "'6.22"
$ws.Range("E20").Value = "  +5.03%  "
$ws.Range("D21").Value = "'70.82"
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("D23").Value = "'234.79"
$ws.Range("E23").Value = "  +1.29%  "
$ws.Range("D24").Value = "'9.72"
$ws.Range("E24").Value = "  +7.33%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  +7.03%  "
$ws.Range("D27").Value = "'2.51"
$ws.Range("E27").Value = "  +12.22%  "
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("D29").Value = "'40.05"
$ws.Range("E29").Value = "  +2.37%  "
$ws.Range("E30").Value = "  +2.92%  "
$ws.Range("D31").Value = "'173.78"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").Value = "'0.0912"
$ws.Range("E32").Value = "  +6.21%  "
$ws.Range("D33").Value = "'21.14"
$ws.Range("E33").Value = "  +3.60%  "
$ws.Range("D34").Value = "'5.69"
$ws.Range("E34").Value = "  +7.03%  "
$ws.Range("E35").Value = "  +1.68%  "
$ws.Range("D36").Value = "'0.112"
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("D37").Value = "'0.0354"
$ws.Range("E37").Value = "  -2.09%  "
$ws.Range("D38").Value = "'4.38"
$ws.Range("E38").Value = "  -2.49%  "
$ws.Range("D39").Value = "'3.53"
$ws.Range("E39").Value = "  +23.14%  "
$ws.Range("D40").Value = "'0.251"
$ws.Range("E40").Value = "  +25.56%  "
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").Value = "'2.19"
$ws.Range("E41").Value = "  +4.50%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "'12.53"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").Value = "'5.46"
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").Value = "'62.19"
$ws.Range("E44").Value = "  -1.57%  "
$ws.Range("E45").Value = "  +4.72%  "
$ws.Range("E46").Value = "  +1.17%  "
$ws.Range("D47").Value = "'100.53"
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("D48").Value = "'1.16"
$ws.Range("E48").Value = "  +4.24%  "
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("D50").Value = "'0.430"
$ws.Range("E50").Value = "  -1.34%  "
$ws.Range("D51").Value = "2.464.08"
$ws.Range("E51").Value = "  +1.94%  "
